$wb = $excel.ActiveWorkbook

# Add a new worksheet at the end of the workbook and name it Pid4CatRecord
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Pid4CatRecord"

$headers = @("landing_page_url", "status", "schema_version", "metadata_license", "curation_contact", "resource_info", "related_identifiers", "change_log")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
